$d = $word.ActiveDocument
$full = $d.WordOpenXML
Write-Output ("doc xml len: " + $full.Length)
Write-Output ("contains lastRenderedPageBreak: " + $full.Contains("lastRenderedPageBreak"))
Write-Output ("contains proofErr: " + $full.Contains("proofErr"))
